# Update "want to go" counts (column F) for two conventions whose numbers
# ticked up between scrapes, on both the "展览" (Exhibition) sheet and the
# "全部类型" (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F4").Value = 609
$wsExhibition.Range("F7").Value = 2729
$wsExhibition.Range("F9").Value = 7619
$wsExhibition.Range("F13").Value = 288

# Sheet 4: 全部类型 (All types) — same events, different row numbers
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value = 609
$wsAll.Range("F9").Value = 2729
$wsAll.Range("F11").Value = 7619
$wsAll.Range("F17").Value = 288
